{"js": "// Apply the \"Added many more features\" edit: retitle the headline/CTA text\n// and swap out the \"What we like\" / \"What we don't like\" bullet copy plus\n// the closing meta-description sentence.\nconst body = context.document.body;\n\nconst replacements = [\n  // Heading1 title AND the later bold CTA repeat the exact same sentence,\n  // both need to become the new title (search() below finds both hits).\n  {\n    from: \"Play De Magorum free and discover its magical features\",\n    to: \"Play De Magorum Slot Game for Free\",\n  },\n  // \"What we like\" bullets\n  {\n    from: \"High RTP of 94.38%\",\n    to: \"Exciting gameplay moments for slot game lovers\",\n  },\n  {\n    from: \"Exciting and rewarding Bonus Game feature\",\n    to: \"Thematic symbols and great payouts\",\n  },\n  {\n    from: \"Wide range of betting options\",\n    to: \"Generous special symbols and multipliers\",\n  },\n  {\n    from: \"Thematic symbols and special features\",\n    to: \"Unique and lucrative Bonus Game feature\",\n  },\n  // \"What we don't like\" bullets\n  {\n    from: \"May not appeal to players who prefer games with higher volatility\",\n    to: \"Limited betting options\",\n  },\n  {\n    from: \"Free spins round can be difficult to trigger\",\n    to: \"Medium volatility may not appeal to all players\",\n  },\n  // Closing italic meta-description sentence\n  {\n    from: \"Read our review of De Magorum, play for free, and experience its magical symbols and Bonus Game feature with a high RTP and wide betting range.\",\n    to: \"Discover the exciting gameplay and lucrative features of De Magorum. Play for free now!\",\n  },\n];\n\nfor (const { from, to } of replacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edit: retitle the headline/CTA text\n# and swap out the \"What we like\" / \"What we don't like\" bullet copy plus\n# the closing meta-description sentence.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    # Heading1 title AND the later bold CTA repeat the exact same sentence;\n    # wdReplaceAll below rewrites both occurrences in one Execute() call.\n    @{ From = \"Play De Magorum free and discover its magical features\"; To = \"Play De Magorum Slot Game for Free\" },\n    # \"What we like\" bullets\n    @{ From = \"High RTP of 94.38%\"; To = \"Exciting gameplay moments for slot game lovers\" },\n    @{ From = \"Exciting and rewarding Bonus Game feature\"; To = \"Thematic symbols and great payouts\" },\n    @{ From = \"Wide range of betting options\"; To = \"Generous special symbols and multipliers\" },\n    @{ From = \"Thematic symbols and special features\"; To = \"Unique and lucrative Bonus Game feature\" },\n    # \"What we don't like\" bullets\n    @{ From = \"May not appeal to players who prefer games with higher volatility\"; To = \"Limited betting options\" },\n    @{ From = \"Free spins round can be difficult to trigger\"; To = \"Medium volatility may not appeal to all players\" },\n    # Closing italic meta-description sentence\n    @{ From = \"Read our review of De Magorum, play for free, and experience its magical symbols and Bonus Game feature with a high RTP and wide betting range.\"; To = \"Discover the exciting gameplay and lucrative features of De Magorum. Play for free now!\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.From\n    $find.Replacement.Text = $r.To\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, \"wdReplaceAll\")\n}\n"}
